$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new profit row produced by the 2025-09-07 run.
# Column A holds the date as literal text (e.g. "08/19/2025" in the
# existing rows), so we force a text number format before assigning the
# value to stop Excel from auto-parsing the string into a date serial.
# ClearFormats() afterwards drops the temporary formatting again so the
# new cell ends up unstyled, just like the other date cells above it.
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "09/07/2025"
$ws.Range("A21").ClearFormats()

# Column B is the plain numeric profit value for that day.
$ws.Range("B21").Value = 14447.36
